$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values must stay text -- several look numeric (e.g. "0.9992",
# "5.300", "23.477.41") and would otherwise be reinterpreted/truncated by Excel's
# type coercion, so each D-cell is forced to Text format before the assignment.
$updates = [ordered]@{
    'D2' = '23.477.41'
    'E2' = '  -0.28%  '
    'D3' = '1.632.49'
    'E3' = '  -0.50%  '
    'D4' = '0.9992'
    'E4' = '  +0.00%  '
    'D5' = '0.9989'
    'E5' = '  -0.04%  '
    'D6' = '304.59'
    'E6' = '  -1.39%  '
    'D7' = '0.3788'
    'E7' = '  +0.33%  '
    'D8' = '52.01'
    'E8' = '  -1.52%  '
    'D9' = '0.3646'
    'E9' = '  -1.10%  '
    'D10' = '1.236'
    'E10' = '  -3.44%  '
    'D11' = '0.08122'
    'E11' = '  -1.17%  '
    'D12' = '0.9991'
    'E12' = '  +0.00%  '
    'D13' = '22.71'
    'E13' = '  -2.44%  '
    'D14' = '6.592'
    'E14' = '  -1.27%  '
    'D15' = '0.00001251'
    'E15' = '  -2.69%  '
    'D16' = '7.273'
    'E16' = '  -2.81%  '
    'D17' = '1.633.85'
    'E17' = '  -0.40%  '
    'D18' = '93.89'
    'E18' = '  -1.27%  '
    'D19' = '0.06943'
    'E19' = '  -0.28%  '
    'E20' = '  -2.65%  '
    'D21' = '6.447'
    'E21' = '  -2.21%  '
    'D22' = '0.9990'
    'E22' = '  +0.09%  '
    'D23' = '23.485.99'
    'E23' = '  -0.22%  '
    'E24' = '  -1.57%  '
    'D25' = '3.241'
    'E25' = '  +4.57%  '
    'D26' = '2.431'
    'E26' = '  +0.50%  '
    'D27' = '21.29'
    'E27' = '  -0.61%  '
    'D28' = '149.51'
    'E28' = '  -1.36%  '
    'D29' = '5.300'
    'E29' = '  -0.46%  '
    'D30' = '134.79'
    'E30' = '  -1.06%  '
    'D31' = '2.305'
    'E31' = '  -4.69%  '
    'D32' = '1.806.53'
    'E32' = '  -0.83%  '
    'D33' = '6.863'
    'E33' = '  -0.12%  '
    'D34' = '11.04'
    'E34' = '  +5.46%  '
    'D35' = '0.9643'
    'E35' = '  -1.98%  '
    'E36' = '  -0.51%  '
    'D37' = '0.2544'
    'E37' = '  -0.25%  '
    'B38' = 'Hedera'
    'C38' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D38' = '0.07223'
    'E38' = '  -3.66%  '
    'B39' = 'Stellar'
    'C39' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D39' = '0.08836'
    'E39' = '  -0.54%  '
    'D40' = '6.127'
    'E40' = '  -1.69%  '
    'D41' = '0.7123'
    'E41' = '  -0.77%  '
    'D42' = '1.357'
    'E42' = '  -3.10%  '
    'D43' = '16.41'
    'E43' = '  +1.41%  '
    'D44' = '12.38'
    'E44' = '  -2.29%  '
    'D45' = '0.6540'
    'E45' = '  -1.45%  '
    'D46' = '2.345'
    'E46' = '  -1.22%  '
    'D47' = '0.9983'
    'E47' = '  +0.00%  '
    'D48' = '4.005'
    'E48' = '  -1.05%  '
    'D49' = '0.08016'
    'E49' = '  -0.72%  '
    'D50' = '1.210'
    'E50' = '  -1.25%  '
    'D51' = '125.87'
    'E51' = '  -4.06%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($ref -match '^D') {
        $cell.NumberFormat = '@'
    }
    $cell.Value = $updates[$ref]
}
